$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also updates the tab name "Through 2022-08-18" -> "Through 2022-08-19")
$ws.Name = "Through 2022-08-19"

# Update the header label in I1 ("2022 (through 08-18)" -> "2022 (through 08-19)")
$ws.Range("I1").Value = "2022 (through 08-19)"

# Update the September value (row 9) in the 2022 column
$ws.Range("I9").Value = 101

# Update the Total value (row 14) in the 2022 column
$ws.Range("I14").Value = 1072
